# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on the two existing sheets to
#    their new, more specific names.
# 2. Add a third sheet "PO Forecast" with a Prophet-style forecast table
#    (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet -------
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# --- Header row (copy the header formatting used on "Weekly Quantity") ---
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$wsForecast.Range("A1").Value = "ds"

$wsWeekly.Range("B1").Copy()
$wsForecast.Range("B1").PasteSpecial(-4122)
$wsForecast.Range("B1").Value = "PO_Forecast"

$wsWeekly.Range("B1").Copy()
$wsForecast.Range("C1").PasteSpecial(-4122)
$wsForecast.Range("C1").Value = "yhat_lower"

$wsWeekly.Range("B1").Copy()
$wsForecast.Range("D1").PasteSpecial(-4122)
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows -------------------------------------------------------------
$data = @(
    @(45312.99999999999, 145, 55.80605007644771,   239.4233186585254),
    @(45319.99999999999, 135, 46.46930903330172,   232.5003062151852),
    @(45354.99999999999, 87,  -10.08428757415571,  176.5280870566307),
    @(45396.99999999999, 30,  -59.03796949461427,  125.9179461726418),
    @(45403.99999999999, 21,  -74.21846571525892,  114.8189333997891),
    @(45410.99999999999, 11,  -82.81038533042869,  107.9686444079915),
    @(45417.99999999999, 1,   -91.94832323495083,  96.84633198818436),
    @(45424.99999999999, 0,   -104.9105311462103,  85.96288164144651),
    @(45431.99999999999, 0,   -116.1165701015045,  74.44036188288177),
    @(45438.99999999999, 0,   -118.9669646951618,  61.48151544848211),
    @(45445.99999999999, 0,   -129.5421049898719,  53.13989258670562),
    @(45452.99999999999, 0,   -133.216081186481,   42.23425427779846),
    @(45459.99999999999, 0,   -145.5770808263345,  32.84228569006386)
)

$row = 2
foreach ($r in $data) {
    # Column A gets the same date/time number-format as the other sheets.
    $wsWeekly.Range("A2").Copy()
    $wsForecast.Range("A$row").PasteSpecial(-4122)
    $wsForecast.Range("A$row").Value = $r[0]

    $wsForecast.Range("B$row").Value = $r[1]
    $wsForecast.Range("C$row").Value = $r[2]
    $wsForecast.Range("D$row").Value = $r[3]

    $row = $row + 1
}
